$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (adds new shared strings at the end, indices 4-7)
$ws.Range("E1").Value = "manufacturing_date"
$ws.Range("F1").Value = "expiry_date"
$ws.Range("G1").Value = "batch_id"
$ws.Range("H1").Value = "unit"

# Center-align the whole header row (A1:D1 + G1 -> style w/o wrap, E1/F1/H1 -> style w/ wrap)
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("G1").HorizontalAlignment = -4108

$ws.Range("E1:F1").HorizontalAlignment = -4108
$ws.Range("E1:F1").WrapText = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").WrapText = $true

# New column widths for the two newly introduced data columns
$ws.Columns("E").ColumnWidth = 18.52
$ws.Columns("F").ColumnWidth = 11.46

# Move the active selection to E2 (below the newly added "manufacturing_date" header)
$ws.Range("E2").Select()
